$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'59.456.16"
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").Value = "'2.638.74"
$ws.Range("E3").Value = "  +1.29%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'536.82"
$ws.Range("E5").Value = "  -0.82%  "
$ws.Range("D6").Value = "'145.20"
$ws.Range("E6").Value = "  +2.47%  "
$ws.Range("E8").Value = "  +1.04%  "
$ws.Range("D9").Value = "'6.99"
$ws.Range("E9").Value = "  +8.45%  "
$ws.Range("E10").Value = "  -1.45%  "
$ws.Range("E11").Value = "  +0.64%  "
$ws.Range("D12").Value = "'0.135"
$ws.Range("E12").Value = "  +0.24%  "
$ws.Range("D13").Value = "'3.108.33"
$ws.Range("E13").Value = "  +1.49%  "
$ws.Range("D14").Value = "'59.381.57"
$ws.Range("E14").Value = "  +0.35%  "
$ws.Range("D15").Value = "'21.31"
$ws.Range("E15").Value = "  +3.28%  "
$ws.Range("D16").Value = "'2.659.02"
$ws.Range("E16").Value = "  +0.61%  "
$ws.Range("E17").Value = "  +0.58%  "
$ws.Range("D18").Value = "'4.50"
$ws.Range("E18").Value = "  +2.63%  "
$ws.Range("D19").Value = "'337.59"
$ws.Range("E19").Value = "  -1.35%  "
$ws.Range("E20").Value = "  +1.43%  "
$ws.Range("D21").Value = "'6.23"
$ws.Range("E21").Value = "  -2.44%  "
$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").Value = "'66.31"
$ws.Range("E23").Value = "  -1.85%  "
$ws.Range("E24").Value = "  +1.96%  "
$ws.Range("E25").Value = "  -0.32%  "
$ws.Range("E26").Value = "  -0.17%  "
$ws.Range("D27").Value = "'7.28"
$ws.Range("E27").Value = "  +0.80%  "
$ws.Range("D28").Value = "'0.0₃0750"
$ws.Range("E28").Value = "  -2.25%  "
$ws.Range("D29").Value = "'0.998"
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("D30").Value = "'1.66"
$ws.Range("E30").Value = "  -3.19%  "
$ws.Range("D31").Value = "'5.90"
$ws.Range("E31").Value = "  +1.15%  "
$ws.Range("D32").Value = "'18.83"
$ws.Range("E32").Value = "  +0.41%  "
$ws.Range("D33").Value = "'151.19"
$ws.Range("E33").Value = "  +1.12%  "
$ws.Range("D34").Value = "'4.01"
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("E35").Value = "  +1.91%  "
$ws.Range("D36").Value = "'0.842"
$ws.Range("E36").Value = "  +2.19%  "
$ws.Range("D37").Value = "'0.836"
$ws.Range("E37").Value = "  -0.40%  "
$ws.Range("E38").Value = "  -1.23%  "
$ws.Range("D39").Value = "'3.61"
$ws.Range("E39").Value = "  +0.97%  "
$ws.Range("D40").Value = "'284.93"
$ws.Range("E40").Value = "  +3.80%  "
$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("E42").Value = "  +0.70%  "
$ws.Range("D43").Value = "'10.75"
$ws.Range("E43").Value = "  +0.09%  "
$ws.Range("D44").Value = "'0.0538"
$ws.Range("E44").Value = "  +2.27%  "
$ws.Range("D45").Value = "'19.19"
$ws.Range("E45").Value = "  +2.39%  "
$ws.Range("E46").Value = "  -1.71%  "
$ws.Range("D47").Value = "'0.0227"
$ws.Range("E47").Value = "  +1.27%  "
$ws.Range("D48").Value = "'1.960.56"
$ws.Range("E48").Value = "  +0.49%  "
$ws.Range("D49").Value = "'4.56"
$ws.Range("E49").Value = "  +0.73%  "
$ws.Range("D50").Value = "'18.39"
$ws.Range("E50").Value = "  -1.10%  "
$ws.Range("D51").Value = "'111.55"
$ws.Range("E51").Value = "  +0.13%  "
